# Restore cell C10 on the "Rules" sheet to its prior value (1), as part of
# reverting the table's min-threshold for the "R30" rule row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = 1
